# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the data block (rows 203-204),
# pushing the existing rows 203:328 down to 205:330.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 203 (existing data shifts down by 2 rows).
$ws.Rows("203:204").Insert()

# --- New row 203 ---
$ws.Cells.Item(203,1).Value  = 3
$ws.Cells.Item(203,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(203,3).Value  = "Coquimbo"
$ws.Cells.Item(203,4).Value  = 44830
$ws.Cells.Item(203,5).Value  = 5
$ws.Cells.Item(203,6).Value  = "Fruta"
$ws.Cells.Item(203,7).Value  = 100101
$ws.Cells.Item(203,8).Value  = "Berries"
$ws.Cells.Item(203,9).Value  = 100112025
$ws.Cells.Item(203,10).Value = "Frutilla"
$ws.Cells.Item(203,11).Value = "Sin especificar"
$ws.Cells.Item(203,12).Value = "Especial"
$ws.Cells.Item(203,13).Value = 56
$ws.Cells.Item(203,14).Value = 14000
$ws.Cells.Item(203,15).Value = 14000
$ws.Cells.Item(203,16).Value = 14000
$ws.Cells.Item(203,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(203,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(203,19).Value = 2000
$ws.Cells.Item(203,20).Value = 7

# --- New row 204 ---
$ws.Cells.Item(204,1).Value  = 3
$ws.Cells.Item(204,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(204,3).Value  = "Coquimbo"
$ws.Cells.Item(204,4).Value  = 44830
$ws.Cells.Item(204,5).Value  = 5
$ws.Cells.Item(204,6).Value  = "Fruta"
$ws.Cells.Item(204,7).Value  = 100101
$ws.Cells.Item(204,8).Value  = "Berries"
$ws.Cells.Item(204,9).Value  = 100112025
$ws.Cells.Item(204,10).Value = "Frutilla"
$ws.Cells.Item(204,11).Value = "Sin especificar"
$ws.Cells.Item(204,12).Value = "Primera"
$ws.Cells.Item(204,13).Value = 60
$ws.Cells.Item(204,14).Value = 12000
$ws.Cells.Item(204,15).Value = 12000
$ws.Cells.Item(204,16).Value = 12000
$ws.Cells.Item(204,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(204,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(204,19).Value = 1714
$ws.Cells.Item(204,20).Value = 7
